$wb = $excel.ActiveWorkbook

# --- Constants sheet: add MaxRetryNumberDispatcher row (new setting, value 2) ---
$wsConstants = $wb.Worksheets.Item("Constants")
$wsConstants.Range("A4").Value = "MaxRetryNumberDispatcher"
$wsConstants.Range("B4").Value = 2
$wsConstants.Range("A5").Select() | Out-Null

# --- Assets sheet: add EmailRecipients + ExScreenshotsFolderPath rows ---
$wsAssets = $wb.Worksheets.Item("Assets")
$wsAssets.Range("A3").Value = "EmailRecipients"
$wsAssets.Range("B3").Value = "EmailRecipients"
$wsAssets.Range("C3").Value = "Academy-REFramework"
$wsAssets.Range("A4").Value = "ExScreenshotsFolderPath"
$wsAssets.Range("B4").Value = "ExScreenshotsFolderPath"
$wsAssets.Range("C4").Value = "Academy-REFramework"
$wsAssets.Range("C4").Select() | Out-Null

# --- Settings sheet becomes the active/selected tab (done last so it "wins" as the active sheet) ---
$wsSettings = $wb.Worksheets.Item("Settings")
$wsSettings.Activate()
$wsSettings.Range("A6").Select() | Out-Null
